$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typos in the "details_of_publiclaw" text shared by rows 12 and 13
# (D12, D13): "suthority" -> "authority", "rasied" -> "raised"
$fixedText = "Located at Title 3, Sec. 301-1A subsection(a); `nif the Treasury exercises the authority to borrow an `nadditional 900B dollars, debt limit is raised by 400B(August 2011) and if disapproval committee has lapsed discusssions debt limit is raised an additional 500B (Septmber 2011) "
$ws.Range("D12").Value = $fixedText
$ws.Range("D13").Value = $fixedText

# Add the missing public_law value "None" for row 28
$ws.Range("C28").Value = "None"

# Update the saved selection to D14 (matches the workbook's last active cell)
$ws.Range("D14").Select()
